{"js": "// Apply spell/grammar-check proofing marks (w:proofErr) around several\n// list items, and split a few runs so the flagged word/phrase sits in\n// its own run \u2014 matching what Word's background proofer inserts when it\n// flags \"Sambata\", \"Duminica\", \"Luni\", \"healthbars\", \"TurretDestroyer\",\n// \"powerpoint\" as misspelled and \"Set  up\", \"360 degree\", \"60 page\" as\n// grammar issues. Also greens-up (color 92D050) the \"TurretDestroyer\" /\n// \"360 degree attack\" sub-steps under Raid logic to mark them done.\n//\n// Word's JS API has no typed object for w:proofErr (it is purely a\n// proofing-UI marker with no selectable range), so the only way to\n// reproduce it faithfully is to splice the paragraph's OOXML directly\n// via Range/Paragraph.insertOoxml(..., Replace).\n\nconst W_NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\nfunction pkg(pXml) {\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' + pXml + '</pkg:xmlData></pkg:part></pkg:package>'\n  );\n}\n\n// Replace the whole paragraph at `items[index]` with the literal OOXML\n// for `<w:p>` given in `innerPPlusRuns` (the full <w:p ...>...</w:p>).\nfunction replaceParagraph(paragraph, pElementXml) {\n  paragraph.insertOoxml(pkg(pElementXml), Word.InsertLocation.replace);\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nfunction findIndex(text, occurrence) {\n  let seen = 0;\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text === text) {\n      if (seen === occurrence) return i;\n      seen++;\n    }\n  }\n  throw new Error(\"paragraph not found: \" + text + \" occurrence \" + occurrence);\n}\n\n// --- 1. \"Sambata\" -> wrap the single run in spellStart/spellEnd ---\n{\n  const i = findIndex(\"Sambata\", 0);\n  replaceParagraph(items[i],\n    `<w:p ${W_NS}>` +\n      `<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>` +\n      `<w:proofErr w:type=\"spellStart\"/>` +\n      `<w:r><w:t>Sambata</w:t></w:r>` +\n      `<w:proofErr w:type=\"spellEnd\"/>` +\n    `</w:p>`);\n}\n\n// --- 2. \"Duminica\" -> wrap the single run in spellStart/spellEnd ---\n{\n  const i = findIndex(\"Duminica\", 0);\n  replaceParagraph(items[i],\n    `<w:p ${W_NS}>` +\n      `<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>` +\n      `<w:proofErr w:type=\"spellStart\"/>` +\n      `<w:r><w:t>Duminica</w:t></w:r>` +\n      `<w:proofErr w:type=\"spellEnd\"/>` +\n    `</w:p>`);\n}\n\n// --- 3. \"Fix healthbars\" -> \"Fix \" + spellStart/\"healthbars\"/spellEnd ---\n{\n  const i = findIndex(\"Fix healthbars\", 0);\n  replaceParagraph(items[i],\n    `<w:p ${W_NS}>` +\n      `<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"1\"/><w:numId w:val=\"1\"/></w:numPr>` +\n      `<w:rPr><w:color w:val=\"92D050\"/></w:rPr></w:pPr>` +\n      `<w:r><w:rPr><w:color w:val=\"92D050\"/></w:rPr><w:t xml:space=\"preserve\">Fix </w:t></w:r>` +\n      `<w:proofErr w:type=\"spellStart\"/>` +\n      `<w:r><w:rPr><w:color w:val=\"92D050\"/></w:rPr><w:t>healthbars</w:t></w:r>` +\n      `<w:proofErr w:type=\"spellEnd\"/>` +\n    `</w:p>`);\n}\n\n// --- 4. \"Luni\" -> wrap the single run in spellStart/spellEnd ---\n{\n  const i = findIndex(\"Luni\", 0);\n  replaceParagraph(items[i],\n    `<w:p ${W_NS}>` +\n      `<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>` +\n      `<w:proofErr w:type=\"spellStart\"/>` +\n      `<w:r><w:t>Luni</w:t></w:r>` +\n      `<w:proofErr w:type=\"spellEnd\"/>` +\n    `</w:p>`);\n}\n\n// --- 5. \"Set  up turret and damage enemies with it\" -> gramStart/\"Set  up\"/gramEnd + \" turret...\" ---\n{\n  const i = findIndex(\"Set  up turret and damage enemies with it\", 0);\n  replaceParagraph(items[i],\n    `<w:p ${W_NS}>` +\n      `<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"2\"/><w:numId w:val=\"1\"/></w:numPr>` +\n      `<w:rPr><w:color w:val=\"92D050\"/></w:rPr></w:pPr>` +\n      `<w:proofErr w:type=\"gramStart\"/>` +\n      `<w:r><w:rPr><w:color w:val=\"92D050\"/></w:rPr><w:t>Set  up</w:t></w:r>` +\n      `<w:proofErr w:type=\"gramEnd\"/>` +\n      `<w:r><w:rPr><w:color w:val=\"92D050\"/></w:rPr><w:t xml:space=\"preserve\"> turret and damage enemies with it</w:t></w:r>` +\n    `</w:p>`);\n}\n\n// --- 6. First \"TurretDestroyer\" (Raid logic, ilvl 3, no color) ---\n//      -> spellStart/spellEnd AND turns green (92D050) same as its sibling\n{\n  const i = findIndex(\"TurretDestroyer\", 0);\n  replaceParagraph(items[i],\n    `<w:p ${W_NS}>` +\n      `<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"3\"/><w:numId w:val=\"1\"/></w:numPr>` +\n      `<w:rPr><w:color w:val=\"92D050\"/></w:rPr></w:pPr>` +\n      `<w:proofErr w:type=\"spellStart\"/>` +\n      `<w:r><w:rPr><w:color w:val=\"92D050\"/></w:rPr><w:t>TurretDestroyer</w:t></w:r>` +\n      `<w:proofErr w:type=\"spellEnd\"/>` +\n    `</w:p>`);\n}\n\n// --- 7. First \"360 degree attack\" (Raid logic, ilvl 3, no color) ---\n//      -> gramStart/\"360 degree\"/gramEnd + \" attack\" AND turns green\n{\n  const i = findIndex(\"360 degree attack\", 0);\n  replaceParagraph(items[i],\n    `<w:p ${W_NS}>` +\n      `<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"3\"/><w:numId w:val=\"1\"/></w:numPr>` +\n      `<w:rPr><w:color w:val=\"92D050\"/></w:rPr></w:pPr>` +\n      `<w:proofErr w:type=\"gramStart\"/>` +\n      `<w:r><w:rPr><w:color w:val=\"92D050\"/></w:rPr><w:t>360 degree</w:t></w:r>` +\n      `<w:proofErr w:type=\"gramEnd\"/>` +\n      `<w:r><w:rPr><w:color w:val=\"92D050\"/></w:rPr><w:t xml:space=\"preserve\"> attack</w:t></w:r>` +\n    `</w:p>`);\n}\n\n// --- 8. Second \"TurretDestroyer\" (Animations, ilvl 4, already green) ---\n//      -> spellStart/spellEnd only\n{\n  const i = findIndex(\"TurretDestroyer\", 1);\n  replaceParagraph(items[i],\n    `<w:p ${W_NS}>` +\n      `<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"4\"/><w:numId w:val=\"1\"/></w:numPr>` +\n      `<w:rPr><w:color w:val=\"92D050\"/></w:rPr></w:pPr>` +\n      `<w:proofErr w:type=\"spellStart\"/>` +\n      `<w:r><w:rPr><w:color w:val=\"92D050\"/></w:rPr><w:t>TurretDestroyer</w:t></w:r>` +\n      `<w:proofErr w:type=\"spellEnd\"/>` +\n    `</w:p>`);\n}\n\n// --- 9. Second \"360 degree attack\" (Animations, ilvl 4, already green) ---\n//      -> gramStart/\"360 degree\"/gramEnd + \" attack\"\n{\n  const i = findIndex(\"360 degree attack\", 1);\n  replaceParagraph(items[i],\n    `<w:p ${W_NS}>` +\n      `<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"4\"/><w:numId w:val=\"1\"/></w:numPr>` +\n      `<w:rPr><w:color w:val=\"92D050\"/></w:rPr></w:pPr>` +\n      `<w:proofErr w:type=\"gramStart\"/>` +\n      `<w:r><w:rPr><w:color w:val=\"92D050\"/></w:rPr><w:t>360 degree</w:t></w:r>` +\n      `<w:proofErr w:type=\"gramEnd\"/>` +\n      `<w:r><w:rPr><w:color w:val=\"92D050\"/></w:rPr><w:t xml:space=\"preserve\"> attack</w:t></w:r>` +\n    `</w:p>`);\n}\n\n// --- 10. \"Write 60 page documentation\" -> \"Write \" + gramStart/\"60 page\"/gramEnd + \" documentation\" ---\n{\n  const i = findIndex(\"Write 60 page documentation\", 0);\n  replaceParagraph(items[i],\n    `<w:p ${W_NS}>` +\n      `<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"1\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>` +\n      `<w:r><w:t xml:space=\"preserve\">Write </w:t></w:r>` +\n      `<w:proofErr w:type=\"gramStart\"/>` +\n      `<w:r><w:t>60 page</w:t></w:r>` +\n      `<w:proofErr w:type=\"gramEnd\"/>` +\n      `<w:r><w:t xml:space=\"preserve\"> documentation</w:t></w:r>` +\n    `</w:p>`);\n}\n\n// --- 11. \"Make powerpoint presentation\" -> \"Make \" + spellStart/\"powerpoint\"/spellEnd + \" presentation\" ---\n{\n  const i = findIndex(\"Make powerpoint presentation\", 0);\n  replaceParagraph(items[i],\n    `<w:p ${W_NS}>` +\n      `<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"1\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>` +\n      `<w:r><w:t xml:space=\"preserve\">Make </w:t></w:r>` +\n      `<w:proofErr w:type=\"spellStart\"/>` +\n      `<w:r><w:t>powerpoint</w:t></w:r>` +\n      `<w:proofErr w:type=\"spellEnd\"/>` +\n      `<w:r><w:t xml:space=\"preserve\"> presentation</w:t></w:r>` +\n    `</w:p>`);\n}\n\nawait context.sync();\n", "ps1": "# Apply spell/grammar-check proofing marks (w:proofErr) around several\n# list items, and split a few runs so the flagged word/phrase sits in\n# its own run -- matching what Word's background proofer inserts when it\n# flags \"Sambata\", \"Duminica\", \"Luni\", \"healthbars\", \"TurretDestroyer\",\n# \"powerpoint\" as misspelled and \"Set  up\", \"360 degree\", \"60 page\" as\n# grammar issues. Also greens-up (color 92D050) the \"TurretDestroyer\" /\n# \"360 degree attack\" sub-steps under Raid logic to mark them done.\n#\n# The Word object model has no exposed type for w:proofErr (it is purely\n# a proofing-UI marker, not a selectable/addressable range), so the only\n# faithful way to reproduce it is to splice the paragraph's WordOpenXML\n# directly via Range.InsertXML(...), which replaces the range content\n# in place with the supplied OOXML package.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParaIndex($doc, $text, $occurrence) {\n    $seen = 0\n    $count = $doc.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $doc.Paragraphs($i)\n        $t = $p.Range.Text\n        $t = $t.TrimEnd([char]13, [char]7)\n        if ($t -eq $text) {\n            if ($seen -eq $occurrence) {\n                return $i\n            }\n            $seen = $seen + 1\n        }\n    }\n    throw \"paragraph not found: $text occurrence $occurrence\"\n}\n\nfunction Set-ParaXml($doc, $index, $pXml) {\n    $pkgXml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' + $pXml + '</pkg:xmlData></pkg:part></pkg:package>'\n    $p = $doc.Paragraphs($index)\n    $null = $p.Range.InsertXML($pkgXml)\n}\n\n$W = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n\n# --- 1. \"Sambata\" -> wrap the single run in spellStart/spellEnd ---\n$i = Find-ParaIndex $d \"Sambata\" 0\n$xml = \"<w:p $W>\" +\n    \"<w:pPr><w:pStyle w:val=`\"ListParagraph`\"/><w:numPr><w:ilvl w:val=`\"0`\"/><w:numId w:val=`\"1`\"/></w:numPr></w:pPr>\" +\n    \"<w:proofErr w:type=`\"spellStart`\"/>\" +\n    \"<w:r><w:t>Sambata</w:t></w:r>\" +\n    \"<w:proofErr w:type=`\"spellEnd`\"/>\" +\n    \"</w:p>\"\nSet-ParaXml $d $i $xml\n\n# --- 2. \"Duminica\" -> wrap the single run in spellStart/spellEnd ---\n$i = Find-ParaIndex $d \"Duminica\" 0\n$xml = \"<w:p $W>\" +\n    \"<w:pPr><w:pStyle w:val=`\"ListParagraph`\"/><w:numPr><w:ilvl w:val=`\"0`\"/><w:numId w:val=`\"1`\"/></w:numPr></w:pPr>\" +\n    \"<w:proofErr w:type=`\"spellStart`\"/>\" +\n    \"<w:r><w:t>Duminica</w:t></w:r>\" +\n    \"<w:proofErr w:type=`\"spellEnd`\"/>\" +\n    \"</w:p>\"\nSet-ParaXml $d $i $xml\n\n# --- 3. \"Fix healthbars\" -> \"Fix \" + spellStart/\"healthbars\"/spellEnd ---\n$i = Find-ParaIndex $d \"Fix healthbars\" 0\n$xml = \"<w:p $W>\" +\n    \"<w:pPr><w:pStyle w:val=`\"ListParagraph`\"/><w:numPr><w:ilvl w:val=`\"1`\"/><w:numId w:val=`\"1`\"/></w:numPr>\" +\n    \"<w:rPr><w:color w:val=`\"92D050`\"/></w:rPr></w:pPr>\" +\n    \"<w:r><w:rPr><w:color w:val=`\"92D050`\"/></w:rPr><w:t xml:space=`\"preserve`\">Fix </w:t></w:r>\" +\n    \"<w:proofErr w:type=`\"spellStart`\"/>\" +\n    \"<w:r><w:rPr><w:color w:val=`\"92D050`\"/></w:rPr><w:t>healthbars</w:t></w:r>\" +\n    \"<w:proofErr w:type=`\"spellEnd`\"/>\" +\n    \"</w:p>\"\nSet-ParaXml $d $i $xml\n\n# --- 4. \"Luni\" -> wrap the single run in spellStart/spellEnd ---\n$i = Find-ParaIndex $d \"Luni\" 0\n$xml = \"<w:p $W>\" +\n    \"<w:pPr><w:pStyle w:val=`\"ListParagraph`\"/><w:numPr><w:ilvl w:val=`\"0`\"/><w:numId w:val=`\"1`\"/></w:numPr></w:pPr>\" +\n    \"<w:proofErr w:type=`\"spellStart`\"/>\" +\n    \"<w:r><w:t>Luni</w:t></w:r>\" +\n    \"<w:proofErr w:type=`\"spellEnd`\"/>\" +\n    \"</w:p>\"\nSet-ParaXml $d $i $xml\n\n# --- 5. \"Set  up turret and damage enemies with it\" -> gramStart/\"Set  up\"/gramEnd + \" turret...\" ---\n$i = Find-ParaIndex $d \"Set  up turret and damage enemies with it\" 0\n$xml = \"<w:p $W>\" +\n    \"<w:pPr><w:pStyle w:val=`\"ListParagraph`\"/><w:numPr><w:ilvl w:val=`\"2`\"/><w:numId w:val=`\"1`\"/></w:numPr>\" +\n    \"<w:rPr><w:color w:val=`\"92D050`\"/></w:rPr></w:pPr>\" +\n    \"<w:proofErr w:type=`\"gramStart`\"/>\" +\n    \"<w:r><w:rPr><w:color w:val=`\"92D050`\"/></w:rPr><w:t>Set  up</w:t></w:r>\" +\n    \"<w:proofErr w:type=`\"gramEnd`\"/>\" +\n    \"<w:r><w:rPr><w:color w:val=`\"92D050`\"/></w:rPr><w:t xml:space=`\"preserve`\"> turret and damage enemies with it</w:t></w:r>\" +\n    \"</w:p>\"\nSet-ParaXml $d $i $xml\n\n# --- 6. First \"TurretDestroyer\" (Raid logic, ilvl 3, no color) ---\n#      -> spellStart/spellEnd AND turns green (92D050) same as its sibling\n$i = Find-ParaIndex $d \"TurretDestroyer\" 0\n$xml = \"<w:p $W>\" +\n    \"<w:pPr><w:pStyle w:val=`\"ListParagraph`\"/><w:numPr><w:ilvl w:val=`\"3`\"/><w:numId w:val=`\"1`\"/></w:numPr>\" +\n    \"<w:rPr><w:color w:val=`\"92D050`\"/></w:rPr></w:pPr>\" +\n    \"<w:proofErr w:type=`\"spellStart`\"/>\" +\n    \"<w:r><w:rPr><w:color w:val=`\"92D050`\"/></w:rPr><w:t>TurretDestroyer</w:t></w:r>\" +\n    \"<w:proofErr w:type=`\"spellEnd`\"/>\" +\n    \"</w:p>\"\nSet-ParaXml $d $i $xml\n\n# --- 7. First \"360 degree attack\" (Raid logic, ilvl 3, no color) ---\n#      -> gramStart/\"360 degree\"/gramEnd + \" attack\" AND turns green\n$i = Find-ParaIndex $d \"360 degree attack\" 0\n$xml = \"<w:p $W>\" +\n    \"<w:pPr><w:pStyle w:val=`\"ListParagraph`\"/><w:numPr><w:ilvl w:val=`\"3`\"/><w:numId w:val=`\"1`\"/></w:numPr>\" +\n    \"<w:rPr><w:color w:val=`\"92D050`\"/></w:rPr></w:pPr>\" +\n    \"<w:proofErr w:type=`\"gramStart`\"/>\" +\n    \"<w:r><w:rPr><w:color w:val=`\"92D050`\"/></w:rPr><w:t>360 degree</w:t></w:r>\" +\n    \"<w:proofErr w:type=`\"gramEnd`\"/>\" +\n    \"<w:r><w:rPr><w:color w:val=`\"92D050`\"/></w:rPr><w:t xml:space=`\"preserve`\"> attack</w:t></w:r>\" +\n    \"</w:p>\"\nSet-ParaXml $d $i $xml\n\n# --- 8. Second \"TurretDestroyer\" (Animations, ilvl 4, already green) ---\n#      -> spellStart/spellEnd only\n$i = Find-ParaIndex $d \"TurretDestroyer\" 1\n$xml = \"<w:p $W>\" +\n    \"<w:pPr><w:pStyle w:val=`\"ListParagraph`\"/><w:numPr><w:ilvl w:val=`\"4`\"/><w:numId w:val=`\"1`\"/></w:numPr>\" +\n    \"<w:rPr><w:color w:val=`\"92D050`\"/></w:rPr></w:pPr>\" +\n    \"<w:proofErr w:type=`\"spellStart`\"/>\" +\n    \"<w:r><w:rPr><w:color w:val=`\"92D050`\"/></w:rPr><w:t>TurretDestroyer</w:t></w:r>\" +\n    \"<w:proofErr w:type=`\"spellEnd`\"/>\" +\n    \"</w:p>\"\nSet-ParaXml $d $i $xml\n\n# --- 9. Second \"360 degree attack\" (Animations, ilvl 4, already green) ---\n#      -> gramStart/\"360 degree\"/gramEnd + \" attack\"\n$i = Find-ParaIndex $d \"360 degree attack\" 1\n$xml = \"<w:p $W>\" +\n    \"<w:pPr><w:pStyle w:val=`\"ListParagraph`\"/><w:numPr><w:ilvl w:val=`\"4`\"/><w:numId w:val=`\"1`\"/></w:numPr>\" +\n    \"<w:rPr><w:color w:val=`\"92D050`\"/></w:rPr></w:pPr>\" +\n    \"<w:proofErr w:type=`\"gramStart`\"/>\" +\n    \"<w:r><w:rPr><w:color w:val=`\"92D050`\"/></w:rPr><w:t>360 degree</w:t></w:r>\" +\n    \"<w:proofErr w:type=`\"gramEnd`\"/>\" +\n    \"<w:r><w:rPr><w:color w:val=`\"92D050`\"/></w:rPr><w:t xml:space=`\"preserve`\"> attack</w:t></w:r>\" +\n    \"</w:p>\"\nSet-ParaXml $d $i $xml\n\n# --- 10. \"Write 60 page documentation\" -> \"Write \" + gramStart/\"60 page\"/gramEnd + \" documentation\" ---\n$i = Find-ParaIndex $d \"Write 60 page documentation\" 0\n$xml = \"<w:p $W>\" +\n    \"<w:pPr><w:pStyle w:val=`\"ListParagraph`\"/><w:numPr><w:ilvl w:val=`\"1`\"/><w:numId w:val=`\"1`\"/></w:numPr></w:pPr>\" +\n    \"<w:r><w:t xml:space=`\"preserve`\">Write </w:t></w:r>\" +\n    \"<w:proofErr w:type=`\"gramStart`\"/>\" +\n    \"<w:r><w:t>60 page</w:t></w:r>\" +\n    \"<w:proofErr w:type=`\"gramEnd`\"/>\" +\n    \"<w:r><w:t xml:space=`\"preserve`\"> documentation</w:t></w:r>\" +\n    \"</w:p>\"\nSet-ParaXml $d $i $xml\n\n# --- 11. \"Make powerpoint presentation\" -> \"Make \" + spellStart/\"powerpoint\"/spellEnd + \" presentation\" ---\n$i = Find-ParaIndex $d \"Make powerpoint presentation\" 0\n$xml = \"<w:p $W>\" +\n    \"<w:pPr><w:pStyle w:val=`\"ListParagraph`\"/><w:numPr><w:ilvl w:val=`\"1`\"/><w:numId w:val=`\"1`\"/></w:numPr></w:pPr>\" +\n    \"<w:r><w:t xml:space=`\"preserve`\">Make </w:t></w:r>\" +\n    \"<w:proofErr w:type=`\"spellStart`\"/>\" +\n    \"<w:r><w:t>powerpoint</w:t></w:r>\" +\n    \"<w:proofErr w:type=`\"spellEnd`\"/>\" +\n    \"<w:r><w:t xml:space=`\"preserve`\"> presentation</w:t></w:r>\" +\n    \"</w:p>\"\nSet-ParaXml $d $i $xml\n"}
